$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D (predicted_proportions) values
$ws.Range("D4").Value = 0.44
$ws.Range("D5").Value = 0.64
$ws.Range("D6").Value = 0.64
$ws.Range("D7").Value = 0.64
$ws.Range("D8").Value = 0.64
$ws.Range("D9").Value = 0.64

$ws.Range("D13").Value = 0.16
$ws.Range("D14").Value = 0.24
$ws.Range("D15").Value = 0.44
$ws.Range("D16").Value = 0.44
$ws.Range("D17").Value = 0.52
$ws.Range("D18").Value = 0.52
$ws.Range("D19").Value = 0.52

# Update selected cell to D24 (reflects the active cell in the diff)
$ws.Range("D24").Select()
